$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# Update the fitted parameter value (h_p_star) in K2
$ws.Range("K2").Value = 0.37140000000000001

# Activate sheet and move the selection to K3 (matches the saved selection state)
$ws.Activate()
$ws.Range("K3").Select()
